$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2, 5 and 6 are re-sorted (descending by Fecha/date), shifting the
# associated Volumen/Precio/Unidad/Origen/Precio-Kg/Kg-unidad data along.
# Resulting row order (by date, descending): row2 <- old row6, row5 <- old row2, row6 <- old row5.

# Row 2 (now date 44344, previously row 6's data)
$ws.Range("D2").Value = 44344
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 13500
$ws.Range("Q2").Value = "$/caja 18 kilos granel"
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 750
$ws.Range("T2").Value = 18

# Row 5 (now date 44334, previously row 2's data)
$ws.Range("D5").Value = 44334
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 12500
$ws.Range("Q5").Value = "$/caja 12 kilos empedrada"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1042
$ws.Range("T5").Value = 12

# Row 6 (now date 44330, previously row 5's data)
$ws.Range("D6").Value = 44330
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15500
$ws.Range("Q6").Value = "$/caja 18 kilos granel"
$ws.Range("R6").Value = "Provincia de Curicó"
$ws.Range("S6").Value = 861
$ws.Range("T6").Value = 18
